$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Bermagui NSW"
$ws.Cells.Item(2, 2).Value = "Bermis Beachside cafe  2/4 Lamont Street, Bermagui NSW"
$ws.Cells.Item(2, 3).Value = "30/12/2020 9:00am-10:00am"
$ws.Cells.Item(2, 4).Value = "Case dined in cafe"

$ws.Cells.Item(3, 1).Value = "Brighton"
$ws.Cells.Item(3, 2).Value = "Sons of Mary Restaurant  14 Spink St, Brighton VIC 3186"
$ws.Cells.Item(3, 3).Value = "24/12/2020 10:00am-11:05am"
$ws.Cells.Item(3, 4).Value = "Case ate at restaurant"

$ws.Cells.Item(4, 1).Value = "Doveton"
$ws.Cells.Item(4, 2).Value = "Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177"
$ws.Cells.Item(4, 3).Value = "26/12/20 6:30pm"
$ws.Cells.Item(4, 4).Value = "Case attended English service"

$ws.Cells.Item(5, 1).Value = "Eden NSW"
$ws.Cells.Item(5, 2).Value = "Great Southern Hotel  158 Imlay Street, Eden NSW"
$ws.Cells.Item(5, 3).Value = "30/12/2020 5:00pm-6:00pm"
$ws.Cells.Item(5, 4).Value = "Case had dinner here"

$ws.Cells.Item(6, 1).Value = "Lakes Entrance"
$ws.Cells.Item(6, 2).Value = "Albert and Co. Cafe - Bellevue Hotel  201 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(6, 3).Value = "29/12/2020 09:30am-10:45am"
$ws.Cells.Item(6, 4).Value = "Case ate in cafe"

$ws.Cells.Item(7, 1).Value = "Lakes Entrance"
$ws.Cells.Item(7, 2).Value = "Esplanade Resort Lakes Entrance - Hotel bar  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(7, 3).Value = "29/12/2020 05:30pm-05:50pm"
$ws.Cells.Item(7, 4).Value = "Case in hotel bar"

$ws.Cells.Item(8, 1).Value = "Lakes Entrance"
$ws.Cells.Item(8, 2).Value = "Esplanade Resort Lakes Entrance - Pool area  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(8, 3).Value = "29/12/2020 03:30pm-04:30pm"
$ws.Cells.Item(8, 4).Value = "Case visited pool area"

$ws.Cells.Item(9, 1).Value = "Lakes Entrance"
$ws.Cells.Item(9, 2).Value = "The Esplanade Resort and Spa  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(9, 3).Value = "29/12/2020 2:30pm-5:50pm"
$ws.Cells.Item(9, 4).Value = "Case attended site"

$ws.Cells.Item(10, 1).Value = "Lakes Entrance"
$ws.Cells.Item(10, 2).Value = "Wyanga Winery  248 Baades Rd, Lakes Entrance VIC 3909"
$ws.Cells.Item(10, 3).Value = "29/12/2020 1:00pm-2:00pm"
$ws.Cells.Item(10, 4).Value = "Case visited vineyard"

$ws.Cells.Item(11, 1).Value = "Leongatha"
$ws.Cells.Item(11, 2).Value = "Coral Fish and Chips 53 Bair St, Leongatha VIC 3953"
$ws.Cells.Item(11, 3).Value = "26/12/2020 5:30pm-7:30pm"
$ws.Cells.Item(11, 4).Value = "Case worked in store"

$ws.Cells.Item(12, 1).Value = "Leongatha"
$ws.Cells.Item(12, 2).Value = "Coral Fish and Chips 53 Bair St, Leongatha VIC 3953"
$ws.Cells.Item(12, 3).Value = "27/12/2020 4:00pm-7:30pm"
$ws.Cells.Item(12, 4).Value = "Case worked in store"

$ws.Cells.Item(13, 1).Value = "Melbourne"
$ws.Cells.Item(13, 2).Value = "European Bier Cafe City  120 Exhibition Street Melbourne VIC 3000"
$ws.Cells.Item(13, 3).Value = "21/12/2020 8:00pm-9:30pm"
$ws.Cells.Item(13, 4).Value = "Case attended cafe"

$ws.Cells.Item(14, 1).Value = "Melbourne"
$ws.Cells.Item(14, 2).Value = "Fonda Mexican Flinders Lane  31 Flinders Lane Melbourne"
$ws.Cells.Item(14, 3).Value = "29/12/2020 6:00pm-7:00pm"
$ws.Cells.Item(14, 4).Value = "Case attended restaurant"

$ws.Cells.Item(15, 1).Value = "Moorabbin"
$ws.Cells.Item(15, 2).Value = "Minnie Miny Mo Cafe  8 Station Street, Moorabbin"
$ws.Cells.Item(15, 3).Value = "24/12/2020 11:00am-12:00pm"
$ws.Cells.Item(15, 4).Value = "Case attended cafe"

$ws.Cells.Item(16, 1).Value = "Oakleigh"
$ws.Cells.Item(16, 2).Value = "Melissa Oakleigh Restaurant  6 Eaton Mall, Oakleigh VIC 3166"
$ws.Cells.Item(16, 3).Value = "28/12/20 7:30pm-8:15pm"
$ws.Cells.Item(16, 4).Value = "Case dined in restaurant"

$ws.Cells.Item(17, 1).Value = "Southbank"
$ws.Cells.Item(17, 2).Value = "Left Bank Melbourne Restaurant and Cocktail Bar  1 Southbank Boulevard, Southbank"
$ws.Cells.Item(17, 3).Value = "25/12/2020 12:00pm-02:30pm"
$ws.Cells.Item(17, 4).Value = "Case attended bar"

$ws.Cells.Item(18, 1).Value = "Southbank"
$ws.Cells.Item(18, 2).Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Cells.Item(18, 3).Value = "23/112/2020 1:00pm-1:30pm"
$ws.Cells.Item(18, 4).Value = "Case attended restaurant"

